$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "A"
$ws.Range("B3").Value = "A"
$ws.Range("B4").Value = "A"
$ws.Range("B5").Value = "A"
$ws.Range("B6").Value = "A"
$ws.Range("B7").Value = "A"
